$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Board config change: Arduino Uno -> Arduino Mega
$ws.Range("A2").Value = "Arduino Mega"

# Kosten angepasst
$ws.Range("B2").Value = 11.99
$ws.Range("B14").Value = 15

# Update selected cell (cosmetic view state) to match F8
$ws.Range("F8").Select()
